# SCH - MCU USART R/Tx Pin-swap  USART R/Tx pin 할당 miss
#
# Applies the Issue-list / ECO-list updates describing the USART1 TX/RX
# pin-swap mistake on the MCU (U1) and its PCB fix.

$wb = $excel.ActiveWorkbook

$wsIssue = $wb.Worksheets.Item("Issue list")
$wsEco   = $wb.Worksheets.Item("ECO list")

# ---------------------------------------------------------------------
# Issue list sheet - row 18 (sheet row 21): new HW issue about the
# USART1 RX/TX pins being swapped on the MCU.
# ---------------------------------------------------------------------
$wsIssue.Activate()

$wsIssue.Range("D21").Value = "HW"
$wsIssue.Range("H21").Value = "MCU의 Pin할당 miss`n - PA9 : USART1_RX -> USART1_TX`n - PA10 : USART1_TX -> USART1_RX"

# Restore the frozen header pane (rows 1-3) and move the view down so
# that row 13 is the first visible row under the freeze, with G25 as
# the active selected cell (mirrors the author scrolling the list down
# after adding the new rows of data below).
$winIssue = $excel.ActiveWindow
$winIssue.FreezePanes = $false
$wsIssue.Range("A4").Select()
$winIssue.FreezePanes = $true
$winIssue.ScrollRow = 13
$wsIssue.Range("G25").Select()

# ---------------------------------------------------------------------
# ECO list sheet - rows 9 (sheet row 12) and 10 (sheet row 13): the
# actual ECO entries describing the pin swap fix for U1 PIN-42(PA9)
# and PIN-43(PA10).
# ---------------------------------------------------------------------
$wsEco.Activate()

# Row 12 (ECO #9) - PIN-42 (PA9) swapped from USART1_RX to USART1_TX
$wsEco.Range("C12").Value = 43107
$wsEco.Range("E12").Value = "U1"
$wsEco.Range("G12").Value = "PIN-42(PA9) - USART1_RX"
$wsEco.Range("H12").Value = "PIN-42(PA9) - USART1_TX"
$wsEco.Range("I12").Value = "USART1 R/Tx Swap : Off-page 할당 miss"
$wsEco.Range("J12").Value = "PCB 수정"

# Row 13 (ECO #10) - PIN-43 (PA10) swapped from USART1_TX to USART1_RX
$wsEco.Range("C13").Value = 43107
$wsEco.Range("E13").Value = "U1"
$wsEco.Range("G13").Value = "PIN-43(PA10) - USART1_TX"
$wsEco.Range("H13").Value = "PIN-43(PA10) - USART1_RX"
$wsEco.Range("I13").Value = "USART1 R/Tx Swap : Off-page 할당 miss"
$wsEco.Range("J13").Value = "PCB 수정"

# F12/F13 ("Pin-Map" column) use the boxed "-" placeholder style that
# is already used a couple of rows up (F9); copy that cell's format
# (border/number format) across before typing the quote-prefixed dash
# so the cells end up visually identical to the rest of the table.
$wsEco.Range("F9").Copy()
$wsEco.Range("F12").PasteSpecial(-4122)
$wsEco.Range("F13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsEco.Range("F12").Value = "'-"
$wsEco.Range("F13").Value = "'-"

# Update the remembered selection on the ECO list sheet.
$wsEco.Range("G19").Select()

# Leave the workbook focused back on the Issue list tab, matching the
# original file's tabSelected state.
$wsIssue.Activate()
